$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells AE1:BP1 were named "Reading_<code>"; rename to "code_<code>"
$readingNames = @("code_101106", "code_110106", "code_205106", "code_208106", "code_205206", "code_208206", "code_205306", "code_208306", "code_303106", "code_305106", "code_403106", "code_405106", "code_103110", "code_110110", "code_303110", "code_305110", "code_403110", "code_405110", "code_103115", "code_110115", "code_303115", "code_305115", "code_403115", "code_405115", "code_103118", "code_110118", "code_303118", "code_305118", "code_403118", "code_405118", "code_101105", "code_110105", "code_303105", "code_305105", "code_103109", "code_110109", "code_303109", "code_305109")
$startCol1 = 31  # column AE
for ($i = 0; $i -lt $readingNames.Length; $i++) {
    $ws.Cells.Item(1, $startCol1 + $i).Value = $readingNames[$i]
}

# Header cells DC1:ER1 were named "c#_p##_#"; rename to "code_c#_p##_#"
$fieldNames = @("code_c6_p11_6", "code_c6_p12_6", "code_c6_p13_6", "code_c6_p14_6", "code_c6_p15_6", "code_c6_p16_6", "code_c6_p17_6", "code_c7_p11_6", "code_c7_p12_6", "code_c7_p13_6", "code_c7_p14_6", "code_c7_p15_6", "code_c7_p16_6", "code_c7_p17_6", "code_c8_p11_6", "code_c8_p12_6", "code_c8_p13_6", "code_c8_p14_6", "code_c8_p15_6", "code_c8_p17_6", "code_c8_p18_6", "code_c6_p11_10", "code_c6_p12_10", "code_c6_p13_10", "code_c6_p14_10", "code_c6_p15_10", "code_c6_p16_10", "code_c6_p17_10", "code_c7_p11_10", "code_c7_p12_10", "code_c7_p13_10", "code_c7_p14_10", "code_c7_p15_10", "code_c7_p16_10", "code_c7_p17_10", "code_c8_p11_10", "code_c8_p12_10", "code_c8_p13_10", "code_c8_p14_10", "code_c8_p15_10", "code_c8_p17_10", "code_c8_p18_10")
$startCol2 = 107  # column DC
for ($i = 0; $i -lt $fieldNames.Length; $i++) {
    $ws.Cells.Item(1, $startCol2 + $i).Value = $fieldNames[$i]
}
